# Final version for the thesis.
# Move the "Population size" / "Number of iterations" labels that used to
# live in A13/A14 up into H2/H3 (next to the corresponding data rows), then
# clear out the now-empty rows 13/14 and move the active selection to H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value2 = $ws.Range("A13").Value2
$ws.Range("H3").Value2 = $ws.Range("A14").Value2

$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()

$ws.Range("H7").Select()
